# Module Overview slide (slide 4): replace the three placeholder
# "Section N" lines with the real agenda items for this module.
#
# Slide.Shapes.Item(1).TextFrame.TextRange.Text = "..." (whole-text
# assignment) would collapse the three paragraphs into one, so each
# paragraph is retargeted individually via TextRange.Characters(start,len),
# which rewrites only the run text in place and leaves the surrounding
# paragraph/run XML (and the other two paragraphs) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$contentShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 6") {
        $contentShape = $candidate
    }
}

$tr = $contentShape.TextFrame.TextRange

$full1 = $tr.Text
$idx1 = $full1.IndexOf("Section 1") + 1
$r1 = $contentShape.TextFrame.TextRange.Characters($idx1, 9)
$r1.Text = "Multi-Device Hybrid Apps Tooling Overview"

$full2 = $contentShape.TextFrame.TextRange.Text
$idx2 = $full2.IndexOf("Section 2") + 1
$r2 = $contentShape.TextFrame.TextRange.Characters($idx2, 9)
$r2.Text = "Mobile testing"

$full3 = $contentShape.TextFrame.TextRange.Text
$idx3 = $full3.IndexOf("Section 3") + 1
$r3 = $contentShape.TextFrame.TextRange.Characters($idx3, 9)
$r3.Text = "Deploying to Azure"

Write-Output ("Module Overview bullets now: " + $contentShape.TextFrame.TextRange.Text)
